# Generate Report for Handoff
# Refresh the "Latest Handoff Datetime" column (E) for rows whose files
# were (re)handed off, on both the zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "2016-03-23 06:28:31"
$zhcn.Range("E5").Value = "2016-03-23 06:28:31"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4").Value = "2016-03-23 06:28:39"
